# Ajout fonction api durée trajet
# Adds a "Coordinates" sheet (city, longitude, latitude) after "Feuil1".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Create the new sheet right after Feuil1 -----------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Coordinates"

# --- Data: City, Longitude (B), Latitude (C) -------------------------------
# Longitude/Latitude are stored as TEXT (not numbers) in the source file.
# Column order per row reproduces the original authoring order (most rows
# were entered B then C, but two rows were entered C then B).
$rows = @(
    @(1,  "Paris",       "B", "2.333333",           "C", "48.866667"),
    @(2,  "Angers",      "B", "-0.5593",             "C", "47.4698"),
    @(3,  "La Rochelle", "B", "-1.150000",           "C", "46.166667"),
    @(4,  "Bordeaux",    "C", "44.833328 ",          "B", "-0.56667"),
    @(5,  "Biarritz",    "B", "-1.56667",            "C", "43.48333"),
    @(6,  "Pau",         "B", "-0.366667",           "C", "43.300000"),
    @(7,  "Toulouse",    "B", "1.433333",            "C", "43.600000"),
    @(8,  "Montpellier", "B", "3.862038",            "C", "43.62505 "),
    @(9,  "Nimes",       "B", "4.35",                "C", "43.833328 "),
    @(10, "Marseille",   "B", "5.400000",            "C", "3.862038"),
    @(11, "Monaco",      "B", "7.424450755119324",   "C", "43.738347784533"),
    @(12, "Toulon",      "C", "43.116669 ",          "B", "5.93333"),
    @(13, "Lyon",        "B", "4.850000",            "C", "45.750000"),
    @(14, "Avignon",     "B", "4.81667",             "C", "43.950001 ")
)

# Force columns B & C to Text format *before* writing, so the numeric-looking
# strings are not reinterpreted as numbers.
$dataRange = $ws2.Range("B1:C14")
$dataRange.NumberFormat = "@"

foreach ($row in $rows) {
    $r = $row[0]
    $ws2.Cells.Item($r, 1).Value = $row[1]
    $col1 = $row[2]; $val1 = $row[3]
    $col2 = $row[4]; $val2 = $row[5]
    $ws2.Range("$col1$r").Value = $val1
    $ws2.Range("$col2$r").Value = $val2
}

# Remove the explicit Text-number-format style again (keeps the values as
# text/shared-strings but drops the now-unneeded numFmt formatting so the
# cells fall back to the default style).
$dataRange.ClearFormats()

# --- Per-row font overrides -------------------------------------------------
# Row 13 (Lyon): Arial, colour #333333
$r13 = $ws2.Range("B13:C13")
$r13.Font.Color = 3355443
$r13.Font.Name  = "Arial"

# Rows 12 & 14 (Toulon / Avignon): Tahoma 9, colour #000000
$r12 = $ws2.Range("B12:C12")
$r12.Font.Color = 0
$r12.Font.Name  = "Tahoma"
$r12.Font.Size  = 9

$r12.Copy()
$ws2.Range("B14:C14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Column widths (best effort) -------------------------------------------
$ws2.Range("B1:B14").ColumnWidth = 17.7109375
$ws2.Range("C1:C14").ColumnWidth = 15.5703125

# --- Sheet view / selection -------------------------------------------------
$ws1.Range("B13").Select()

$ws2.Activate()
$ws2.Range("K28").Select()
